$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.1278195408837234
$ws.Range("C2").Value = -1.400444400388857
$ws.Range("D2").Value = 0.08695924927729271
$ws.Range("E2").Value = -0.207196275323166
$ws.Range("F2").Value = 0.3329321255475683
$ws.Range("G2").Value = 0.101692144383436
$ws.Range("H2").Value = 0.2070326302046627
$ws.Range("I2").Value = 0.1094798668430039
$ws.Range("J2").Value = 0.6123535359610832
$ws.Range("K2").Value = 0.8706252039810712
$ws.Range("B3").Value = 0.5447147469754596
$ws.Range("C3").Value = 0.001955153691105105
$ws.Range("D3").Value = 0.4339574590591693
$ws.Range("E3").Value = 0.1559537485280918
$ws.Range("F3").Value = 0.2410754634991643
$ws.Range("G3").Value = 0.1347863954460349
$ws.Range("H3").Value = 0.6338877181444649
$ws.Range("I3").Value = 0.8905316652766015
$ws.Range("J3").Value = -0.8342693007396872
$ws.Range("K3").Value = -0.1809867134402706
$ws.Range("B4").Value = 0.4178261858414339
$ws.Range("C4").Value = 0.1320128381891874
$ws.Range("D4").Value = 0.2167499100690802
$ws.Range("E4").Value = 0.1128318376506894
$ws.Range("F4").Value = 0.6140076662505021
$ws.Range("G4").Value = 0.8719046650920586
$ws.Range("H4").Value = -0.8522446856725341
$ws.Range("I4").Value = -0.1986482585427252
$ws.Range("J4").Value = 0.1648850825765775
$ws.Range("K4").Value = -0.7358932949942634
$ws.Range("B5").Value = 0.2987028809941715
$ws.Range("C5").Value = 0.1220291778826379
$ws.Range("D5").Value = 0.5943001037495309
$ws.Range("E5").Value = 0.8426663495950357
$ws.Range("F5").Value = -0.8843890704881969
$ws.Range("G5").Value = -0.2316442249927503
$ws.Range("H5").Value = 0.131645037661246
$ws.Range("I5").Value = -0.7692023903738299
$ws.Range("J5").Value = 0.3011532043707359
$ws.Range("K5").Value = -0.2759122492816765
$ws.Range("B6").Value = 0.5654988129866392
$ws.Range("C6").Value = 0.8405425659164625
$ws.Range("D6").Value = -0.8771061057057743
$ws.Range("E6").Value = -0.2212277120495061
$ws.Range("F6").Value = 0.1430861830955678
$ws.Range("G6").Value = -0.7574275902997809
$ws.Range("H6").Value = 0.3130368251941282
$ws.Range("I6").Value = -0.26399300444011
$ws.Range("J6").Value = -0.4695888087369317
$ws.Range("K6").Value = 0.509559616829137
$ws.Range("B7").Value = -0.9989882668015876
$ws.Range("C7").Value = -0.31307288098283
$ws.Range("D7").Value = 0.06321343652415951
$ws.Range("E7").Value = -0.8323569178484489
$ws.Range("F7").Value = 0.240154187874971
$ws.Range("G7").Value = -0.3360291699081773
$ws.Range("H7").Value = -0.54127534092662
$ws.Range("I7").Value = 0.4380173397222749
$ws.Range("J7").Value = -0.2982683749317745
$ws.Range("K7").Value = -0.3162133022174854
$ws.Range("B8").Value = 0.1585282081815222
$ws.Range("C8").Value = -0.7540892131611641
$ws.Range("D8").Value = 0.3107446784608449
$ws.Range("E8").Value = -0.2689678476130288
$ws.Range("F8").Value = -0.4758419697301518
$ws.Range("G8").Value = 0.5026967102291701
$ws.Range("H8").Value = -0.2339395072814273
$ws.Range("I8").Value = -0.2520479114279362
$ws.Range("J8").Value = -0.2287093058339515
$ws.Range("K8").Value = -0.5939154206293106
$ws.Range("B9").Value = 0.5708672202635848
$ws.Range("C9").Value = -0.1280619160306041
$ws.Range("D9").Value = -0.390397430094944
$ws.Range("E9").Value = 0.5623484475265214
$ws.Range("F9").Value = -0.1862827769703213
$ws.Range("G9").Value = -0.2099693877891357
$ws.Range("H9").Value = -0.1892248457786991
$ws.Range("I9").Value = -0.5556372688865794
$ws.Range("J9").Value = 0.09960614818241281
$ws.Range("K9").Value = -0.03153282219454001
$ws.Range("B10").Value = -0.5880056801767868
$ws.Range("C10").Value = 0.4383886238637034
$ws.Range("D10").Value = -0.2776623976024654
$ws.Range("E10").Value = -0.2869276062278247
$ws.Range("F10").Value = -0.2597897169002801
$ws.Range("G10").Value = -0.623363232311934
$ws.Range("H10").Value = 0.03314284490197711
$ws.Range("I10").Value = -0.09743360247869032
$ws.Range("J10").Value = -0.2068251325899858
$ws.Range("K10").Value = -1.036548768856761
$ws.Range("B11").Value = -0.2102603443096386
$ws.Range("C11").Value = -0.2120499133054097
$ws.Range("D11").Value = -0.1816236201806256
$ws.Range("E11").Value = -0.5437514435559043
$ws.Range("F11").Value = 0.1133900045823698
$ws.Range("G11").Value = -0.01690728797492691
$ws.Range("H11").Value = -0.1261762057981929
$ws.Range("I11").Value = -0.9558460028783662
$ws.Range("J11").Value = -0.3752733566469112
$ws.Range("K11").Value = -0.4586391267071986
$ws.Range("B12").Value = -0.1939803210945991
$ws.Range("C12").Value = -0.5579424334817151
$ws.Range("D12").Value = 0.09844375035972808
$ws.Range("E12").Value = -0.03216374626862839
$ws.Range("F12").Value = -0.1415591146121002
$ws.Range("G12").Value = -0.9712799915822796
$ws.Range("H12").Value = -0.3907277485715169
$ws.Range("I12").Value = -0.4741015533124014
$ws.Range("J12").Value = -0.692005714912784
$ws.Range("K12").Value = 0.8717747456957818
$ws.Range("B13").Value = 0.2707464529491325
$ws.Range("C13").Value = 0.06232838649038364
$ws.Range("D13").Value = -0.08323367572495516
$ws.Range("E13").Value = -0.9297582624152663
$ws.Range("F13").Value = -0.357013373929053
$ws.Range("G13").Value = -0.4440146383915445
$ws.Range("H13").Value = -0.6636041934262084
$ws.Range("I13").Value = 0.8993931982447689
$ws.Range("J13").Value = -0.4369531029409762
$ws.Range("K13").Value = 1.003492844607284
$ws.Range("B14").Value = -0.1777422226074639
$ws.Range("C14").Value = -0.9662726614299468
$ws.Range("D14").Value = -0.3668156737460341
$ws.Range("E14").Value = -0.441522535107442
$ws.Range("F14").Value = -0.6554532109500923
$ws.Range("G14").Value = 0.9101489502747081
$ws.Range("H14").Value = -0.4249983366195054
$ws.Range("I14").Value = 1.015999555749994
$ws.Range("J14").Value = -0.08426874508913906
$ws.Range("K14").Value = -0.294462537504619
$ws.Range("B15").Value = -0.2244109514860759
$ws.Range("C15").Value = -0.3609463815409155
$ws.Range("D15").Value = -0.6024388739472215
$ws.Range("E15").Value = 0.9508768341145255
$ws.Range("F15").Value = -0.3897482266356107
$ws.Range("G15").Value = 1.048807121260573
$ws.Range("H15").Value = -0.05255046847064687
$ws.Range("I15").Value = -0.2632301165777211
$ws.Range("J15").Value = 0.3351343036404092
$ws.Range("K15").Value = 0.009385904711892601
$ws.Range("B16").Value = -0.4245064810212553
$ws.Range("C16").Value = 1.037734429689356
$ws.Range("D16").Value = -0.3468225008504117
$ws.Range("E16").Value = 1.070715795931569
$ws.Range("F16").Value = -0.04069625057600418
$ws.Range("G16").Value = -0.2561834828541579
$ws.Range("H16").Value = 0.3398831621955681
$ws.Range("I16").Value = 0.01303698643746765
$ws.Range("J16").Value = 0.7296571284789803
$ws.Range("K16").Value = 2.425263408987349
$ws.Range("B17").Value = 1.048302449822014
$ws.Range("C17").Value = -0.3380477122786401
$ws.Range("D17").Value = 1.078518199975304
$ws.Range("E17").Value = -0.03338618974149576
$ws.Range("F17").Value = -0.2491226352005299
$ws.Range("G17").Value = 0.3468183386966098
$ws.Range("H17").Value = 0.01990898967064902
$ws.Range("I17").Value = 0.7364974651863733
$ws.Range("J17").Value = 2.432087912865682
$ws.Range("K17").Value = 9.276200387606567
$ws.Range("B18").Value = -0.4450221160943206
$ws.Range("C18").Value = 1.013258577037035
$ws.Range("D18").Value = -0.07897678827447507
$ws.Range("E18").Value = -0.2854533117167836
$ws.Range("F18").Value = 0.3148482439627
$ws.Range("G18").Value = -0.01000723788286506
$ws.Range("H18").Value = 0.7075488269872576
$ws.Range("I18").Value = 2.403595206618287
$ws.Range("J18").Value = 9.247922562641504
$ws.Range("K18").Value = -8.302694243367595
$ws.Range("B19").Value = 1.047131774781531
$ws.Range("C19").Value = -0.0762236280176779
$ws.Range("D19").Value = -0.2970422425675381
$ws.Range("E19").Value = 0.2967268299057859
$ws.Range("F19").Value = -0.03110735124562009
$ws.Range("G19").Value = 0.6850904903007842
$ws.Range("H19").Value = 2.380517493453967
$ws.Range("I19").Value = 9.2245623784926
$ws.Range("J19").Value = -8.326183261142877
$ws.Range("K19").Value = -0.6168281780219288
$ws.Range("B20").Value = -0.333499057523928
$ws.Range("C20").Value = -0.4283215140497347
$ws.Range("D20").Value = 0.2233711486320382
$ws.Range("E20").Value = -0.07834862591279368
$ws.Range("F20").Value = 0.6496436059220461
$ws.Range("G20").Value = 2.35039660686013
$ws.Range("H20").Value = 9.19684659447498
$ws.Range("I20").Value = -8.352812953902838
$ws.Range("J20").Value = -0.642967415097735
$ws.Range("K20").Value = 0.9401478997058104
$ws.Range("B21").Value = -0.4715856280184577
$ws.Range("C21").Value = 0.1950297526051633
$ws.Range("D21").Value = -0.1030246514249161
$ws.Range("E21").Value = 0.6253569729452143
$ws.Range("F21").Value = 2.325773310604686
$ws.Range("G21").Value = 9.171863421534285
$ws.Range("H21").Value = -8.378042913433744
$ws.Range("I21").Value = -0.6683429869094382
$ws.Range("J21").Value = 0.9146927196727432
$ws.Range("K21").Value = -2.051544358723237
$ws.Range("B22").Value = 0.3052209078313319
$ws.Range("C22").Value = -0.04690201748056927
$ws.Range("D22").Value = 0.6594387888791253
$ws.Range("E22").Value = 2.350872694477139
$ws.Range("F22").Value = 9.193302072250027
$ws.Range("G22").Value = -8.358096204043417
$ws.Range("H22").Value = -0.6490043354755446
$ws.Range("I22").Value = 0.9337835445921647
$ws.Range("J22").Value = -2.032554542815784
$ws.Range("K22").Value = 0.2314556895223195
$ws.Range("B23").Value = -0.2037521027456684
$ws.Range("C23").Value = 0.5642710581134855
$ws.Range("D23").Value = 2.296445563811207
$ws.Range("E23").Value = 9.156748084838261
$ws.Range("F23").Value = -8.38576139252662
$ws.Range("G23").Value = -0.6724114751474837
$ws.Range("H23").Value = 0.9124466247161374
$ws.Range("I23").Value = -2.052887140560817
$ws.Range("J23").Value = 0.2116117104903729
$ws.Range("K23").Value = -0.2393130690951168
$ws.Range("B24").Value = 0.5154074445011286
$ws.Range("C24").Value = 2.268676713534536
$ws.Range("D24").Value = 9.13952823980077
$ws.Range("E24").Value = -8.397856508480213
$ws.Range("F24").Value = -0.6819956026315825
$ws.Range("G24").Value = 0.9040896938991861
$ws.Range("H24").Value = -2.060643857958923
$ws.Range("I24").Value = 0.2041484874144732
$ws.Range("J24").Value = -0.2466327700103889
$ws.Range("K24").Value = -0.3868141421360822
